# Apply the crypto-price-table refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Forces the literal string into the cell even when it looks like a number
    # (e.g. "351.76", "1.00"), matching the source data which stores these as text.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "52.118.40"
$ws.Range("D3").Value = "2.894.01"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "351.76"
$ws.Range("E5").Value = "  -0.95%  "
Set-TextValue "D6" "111.47"
$ws.Range("E6").Value = "  +1.75%  "
Set-TextValue "D7" "0.556"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  +0.07%  "
Set-TextValue "D9" "0.624"
$ws.Range("E9").Value = "  +0.20%  "
Set-TextValue "D10" "39.91"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  +0.48%  "
Set-TextValue "D12" "0.0856"
$ws.Range("E12").Value = "  +2.10%  "
Set-TextValue "D13" "19.95"
$ws.Range("E13").Value = "  -0.25%  "
Set-TextValue "D14" "7.77"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "3.348.68"
$ws.Range("E15").Value = "  +3.42%  "
Set-TextValue "D16" "1.00"
$ws.Range("E16").Value = "  +6.50%  "
$ws.Range("D17").Value = "2.884.55"
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("D18").Value = "52.128.46"
$ws.Range("E18").Value = "  +0.95%  "
Set-TextValue "D19" "7.70"
$ws.Range("E19").Value = "  -0.74%  "
Set-TextValue "D20" "3.33"
$ws.Range("E20").Value = "  +4.99%  "
Set-TextValue "D21" "14.50"
$ws.Range("E21").Value = "  +8.45%  "
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  +0.84%  "
Set-TextValue "D23" "70.76"
$ws.Range("E23").Value = "  +0.32%  "
Set-TextValue "D24" "269.53"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E25").Value = "  +0.33%  "
Set-TextValue "D26" "26.54"
$ws.Range("E26").Value = "  +1.67%  "
Set-TextValue "D27" "0.999"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +1.48%  "
Set-TextValue "D30" "38.23"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "6.48"
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D32" "2.24"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +8.69%  "
$ws.Range("E34").Value = "  +10.63%  "
Set-TextValue "D35" "52.96"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("E37").Value = "  -0.14%  "
Set-TextValue "D38" "3.30"
$ws.Range("E38").Value = "  +4.91%  "
Set-TextValue "D39" "18.62"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("E40").Value = "  +2.47%  "
Set-TextValue "D41" "2.66"
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("E42").Value = "  +1.55%  "
Set-TextValue "D43" "22.59"
$ws.Range("E43").Value = "  +3.28%  "
Set-TextValue "D44" "122.19"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("E46").Value = "  +3.96%  "
$ws.Range("D47").Value = "2.200.05"
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("E48").Value = "  +6.08%  "
Set-TextValue "D49" "0.273"
$ws.Range("E49").Value = "  +23.06%  "
Set-TextValue "D50" "0.946"
$ws.Range("E50").Value = "  +2.53%  "
Set-TextValue "D51" "0.0325"
$ws.Range("E51").Value = "  +11.04%  "
